$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (S.no=2, Name=Jhansi, Mobile=8874543985)
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Jhansi"
$ws.Cells.Item(3, 3).Value = 8874543985

# Match the saved selection state from the diff
$ws.Range("C8").Select()
